# Update cryptocurrency price values in column D to match the refreshed
# symbol-list snapshot. Values are numeric-looking strings that must be
# stored as literal text (preserving exact trailing/leading zero digits),
# so we use the classic apostrophe text-prefix to force Excel to keep them
# as text instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.90"
$ws.Range("D3").Value = "'22.53"
$ws.Range("D4").Value = "'5.390"
$ws.Range("D5").Value = "'0.05761"
$ws.Range("D6").Value = "'3.433"
$ws.Range("D7").Value = "'6.327"
$ws.Range("D8").Value = "'0.8100"
$ws.Range("D9").Value = "'0.8954"
$ws.Range("D10").Value = "'0.1443"
$ws.Range("D11").Value = "'0.07373"
$ws.Range("D12").Value = "'0.03130"
$ws.Range("D13").Value = "'0.02989"
$ws.Range("D14").Value = "'0.09413"
$ws.Range("D15").Value = "'3.946"
$ws.Range("D16").Value = "'0.001587"
$ws.Range("D17").Value = "'0.04796"
$ws.Range("D18").Value = "'0.0005851"
$ws.Range("D19").Value = "'0.006351"
$ws.Range("D20").Value = "'0.004067"
$ws.Range("D21").Value = "'0.0009932"
$ws.Range("D24").Value = "'2.197"
$ws.Range("D25").Value = "'0.3273"
$ws.Range("D27").Value = "'0.0004651"
$ws.Range("D40").Value = "'0.03902"
$ws.Range("D41").Value = "'0.006786"
$ws.Range("D43").Value = "'0.002410"
$ws.Range("D44").Value = "'0.006798"
$ws.Range("D45").Value = "'0.00005647"
$ws.Range("D47").Value = "'0.3801"
$ws.Range("D48").Value = "'0.1635"
$ws.Range("D49").Value = "'0.00002100"
